$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.891.75'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '2.522.67'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = '2.524.44'
$ws.Range("E9").Value = '  +1.21%  '
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("E11").Value = '  -2.63%  '
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = '2.963.42'
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '58.913.79'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '2.512.60'
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E23").Value = '  +2.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.419'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").Value = '0.0₃0768'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.64%  '
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.06%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  +1.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.37'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.07'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.811'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '285.09'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.89%  '
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '129.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0504'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0219'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.62%  '
